$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2005.5
$ws.Range("I17").Value = 1498.3334
$ws.Range("J17").Value = 2174.5557
$ws.Range("K17").Value = 4495.0002
$ws.Range("L17").Value = 6523.6671
$ws.Range("M17").Value = -4327.0002
$ws.Range("N17").Value = -6859.6671
$ws.Range("H41").Value = 980.7222
$ws.Range("I41").Value = 911.2
$ws.Range("J41").Value = 1007.46155
$ws.Range("K41").Value = 911.2
$ws.Range("L41").Value = 1007.46155
$ws.Range("M41").Value = -471.2
$ws.Range("N41").Value = -1887.46155
$ws.Range("H64").Value = 4822.8887
$ws.Range("I64").Value = 4656.1816
$ws.Range("J64").Value = 4937.5
$ws.Range("K64").Value = 4656.1816
$ws.Range("L64").Value = 4937.5
$ws.Range("M64").Value = -4408.1816
$ws.Range("N64").Value = -5433.5
$ws.Range("H67").Value = 4822.8887
$ws.Range("I67").Value = 4656.1816
$ws.Range("J67").Value = 4937.5
$ws.Range("K67").Value = 4656.1816
$ws.Range("L67").Value = 4937.5
$ws.Range("M67").Value = -3798.1816
$ws.Range("N67").Value = -6653.5
$ws.Range("H135").Value = 510.25
$ws.Range("I135").Value = 510.25
$ws.Range("K135").Value = 4592.25
$ws.Range("M135").Value = -2057.25
$ws.Range("H137").Value = 1212529.6
$ws.Range("I137").Value = 4410.75
$ws.Range("J137").Value = 1816589.1
$ws.Range("K137").Value = 13232.25
$ws.Range("L137").Value = 5449767.300000001
$ws.Range("M137").Value = -10682.25
$ws.Range("N137").Value = -5454867.300000001
$ws.Range("H138").Value = 1903.7188
$ws.Range("I138").Value = 1470.0416
$ws.Range("J138").Value = 3204.75
$ws.Range("K138").Value = 4410.1248
$ws.Range("L138").Value = 9614.25
$ws.Range("M138").Value = 729.8752000000004
$ws.Range("N138").Value = -19894.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3709.21
$ws.Range("I32").Value = 3162.0674
$ws.Range("J32").Value = 8136.091
$ws.Range("K32").Value = 3162.0674
$ws.Range("L32").Value = 8136.091
$ws.Range("M32").Value = -2875.0674
$ws.Range("N32").Value = -8710.091
$ws.Range("H74").Value = 37996.855
$ws.Range("I74").Value = 68312.734
$ws.Range("J74").Value = 3017
$ws.Range("K74").Value = 68312.734
$ws.Range("L74").Value = 3017
$ws.Range("M74").Value = -67438.734
$ws.Range("N74").Value = -4765
$ws.Range("H77").Value = 37996.855
$ws.Range("I77").Value = 68312.734
$ws.Range("J77").Value = 3017
$ws.Range("K77").Value = 341563.67
$ws.Range("L77").Value = 15085
$ws.Range("M77").Value = -337195.67
$ws.Range("N77").Value = -23821
$ws.Range("H132").Value = 2740.0476
$ws.Range("I132").Value = 1688.9166
$ws.Range("J132").Value = 4141.5557
$ws.Range("K132").Value = 5066.7498
$ws.Range("L132").Value = 12424.6671
$ws.Range("M132").Value = -2536.7498
$ws.Range("N132").Value = -17484.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3279.6667
$ws.Range("I107").Value = 2636.8333
$ws.Range("J107").Value = 5851
$ws.Range("K107").Value = 2636.8333
$ws.Range("L107").Value = 5851
$ws.Range("M107").Value = -716.8332999999998
$ws.Range("N107").Value = -9691
$ws.Range("H134").Value = 1558.9215
$ws.Range("I134").Value = 841.60974
$ws.Range("K134").Value = 2524.82922
$ws.Range("M134").Value = 10.17077999999992

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3226.2727
$ws.Range("I31").Value = 1960
$ws.Range("J31").Value = 5055.3335
$ws.Range("K31").Value = 1960
$ws.Range("L31").Value = 5055.3335
$ws.Range("M31").Value = -1665
$ws.Range("N31").Value = -5645.3335
$ws.Range("H34").Value = 3226.2727
$ws.Range("I34").Value = 1960
$ws.Range("J34").Value = 5055.3335
$ws.Range("K34").Value = 1960
$ws.Range("L34").Value = 5055.3335
$ws.Range("M34").Value = -1758
$ws.Range("N34").Value = -5459.3335
$ws.Range("H86").Value = 74788
$ws.Range("I86").Value = 94201.91
$ws.Range("K86").Value = 94201.91
$ws.Range("M86").Value = -93078.91
$ws.Range("H89").Value = 74788
$ws.Range("I89").Value = 94201.91
$ws.Range("K89").Value = 471009.55
$ws.Range("M89").Value = -465393.55

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 11396.75
$ws.Range("I43").Value = 3529
$ws.Range("K43").Value = 3529
$ws.Range("M43").Value = -3378
$ws.Range("H46").Value = 33333.332
$ws.Range("J46").Value = 33333.332
$ws.Range("L46").Value = 33333.332
$ws.Range("N46").Value = -33645.332
$ws.Range("H57").Value = 19874.75
$ws.Range("I57").Value = 9999.5
$ws.Range("J57").Value = 29750
$ws.Range("K57").Value = 9999.5
$ws.Range("L57").Value = 29750
$ws.Range("M57").Value = -9179.5
$ws.Range("N57").Value = -31390

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1178.1177
$ws.Range("I16").Value = 1218.5333
$ws.Range("J16").Value = 875
$ws.Range("K16").Value = 1218.5333
$ws.Range("L16").Value = 875
$ws.Range("M16").Value = -1048.5333
$ws.Range("N16").Value = -1215
$ws.Range("H22").Value = 7112.4375
$ws.Range("J22").Value = 8653.846
$ws.Range("L22").Value = 8653.846
$ws.Range("N22").Value = -9243.846
$ws.Range("H27").Value = 7112.4375
$ws.Range("J27").Value = 8653.846
$ws.Range("L27").Value = 8653.846
$ws.Range("N27").Value = -8867.846
$ws.Range("H46").Value = 7718.3687
$ws.Range("I46").Value = 14300.125
$ws.Range("J46").Value = 2931.6365
$ws.Range("K46").Value = 14300.125
$ws.Range("L46").Value = 2931.6365
$ws.Range("M46").Value = -14112.125
$ws.Range("N46").Value = -3307.6365
$ws.Range("H55").Value = 4443.4326
$ws.Range("I55").Value = 1233.619
$ws.Range("J55").Value = 8656.3125
$ws.Range("K55").Value = 1233.619
$ws.Range("L55").Value = 8656.3125
$ws.Range("M55").Value = -1060.619
$ws.Range("N55").Value = -9002.3125
$ws.Range("H100").Value = 10714.214
$ws.Range("I100").Value = 11945.363
$ws.Range("K100").Value = 11945.363
$ws.Range("M100").Value = -11404.363

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 76167.336
$ws.Range("J15").Value = 76167.336
$ws.Range("L15").Value = 76167.336
$ws.Range("N15").Value = -76743.336
$ws.Range("H21").Value = 36005.332
$ws.Range("I21").Value = 38999.5
$ws.Range("J21").Value = 30017
$ws.Range("K21").Value = 38999.5
$ws.Range("L21").Value = 30017
$ws.Range("M21").Value = -38764.5
$ws.Range("N21").Value = -30487
$ws.Range("H35").Value = 36005.332
$ws.Range("I35").Value = 38999.5
$ws.Range("J35").Value = 30017
$ws.Range("K35").Value = 38999.5
$ws.Range("L35").Value = 30017
$ws.Range("M35").Value = -38709.5
$ws.Range("N35").Value = -30597
$ws.Range("H54").Value = 29674.75
$ws.Range("I54").Value = 39700
$ws.Range("J54").Value = 26333
$ws.Range("K54").Value = 39700
$ws.Range("L54").Value = 26333
$ws.Range("M54").Value = -39180
$ws.Range("N54").Value = -27373
$ws.Range("H140").Value = 63750
$ws.Range("J140").Value = 59333.332
$ws.Range("L140").Value = 59333.332
$ws.Range("N140").Value = -69693.33199999999
